$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
# no-op test
